$wb = $excel.ActiveWorkbook
$report = $wb.Worksheets.Item("Report")
$quantities = $wb.Worksheets.Item("Quantities")

# Bump the "TEST #" numbering in column A of the Report sheet by one.
# Only the literal seed cells need updating -- the rest of the column is
# driven by formulas (1+A2, A16+0.01, shared formulas, ...) that recompute
# automatically once the seeds change.
$report.Range("A2").Value = 1
$report.Range("A16").Value = 15.01
$report.Range("A30").Value = 16.01
$report.Range("A33").Value = 17.01
$report.Range("A47").Value = 18.01
$report.Range("A62").Value = 19.01
$report.Range("A76").Value = 20.01
$report.Range("A83").Value = 21.01

# Switch the active/selected sheet from Quantities back to Report, and
# update the Report sheet's scroll/selection state.
$report.Activate()
$report.Range("A84").Select()
